# Updates the "想去人数" (want-to-go count) column F across all four sheets
# to reflect a refreshed scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 443
$ws1.Range("F5").Value = 8681
$ws1.Range("F6").Value = 14
$ws1.Range("F7").Value = 11009
$ws1.Range("F10").Value = 20
$ws1.Range("F15").Value = 295
$ws1.Range("F18").Value = 83
$ws1.Range("F20").Value = 418
$ws1.Range("F21").Value = 37
$ws1.Range("F22").Value = 1865
$ws1.Range("F23").Value = 696
$ws1.Range("F24").Value = 609
$ws1.Range("F25").Value = 350
$ws1.Range("F27").Value = 74
$ws1.Range("F30").Value = 1253
$ws1.Range("F31").Value = 23
$ws1.Range("F33").Value = 5
$ws1.Range("F35").Value = 1421
$ws1.Range("F37").Value = 351
$ws1.Range("F38").Value = 293
$ws1.Range("F39").Value = 30
$ws1.Range("F40").Value = 139
$ws1.Range("F41").Value = 531
$ws1.Range("F42").Value = 367
$ws1.Range("F43").Value = 110
$ws1.Range("F44").Value = 811
$ws1.Range("F45").Value = 650
$ws1.Range("F47").Value = 142
$ws1.Range("F48").Value = 128

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F5").Value = 1
$ws2.Range("F14").Value = 24

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 215
$ws3.Range("F4").Value = 347

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 215
$ws4.Range("F5").Value = 347
$ws4.Range("F8").Value = 443
$ws4.Range("F9").Value = 8681
$ws4.Range("F10").Value = 14
$ws4.Range("F11").Value = 11009
$ws4.Range("F13").Value = 20
$ws4.Range("F15").Value = 295
$ws4.Range("F17").Value = 83
$ws4.Range("F19").Value = 418
$ws4.Range("F20").Value = 1865
$ws4.Range("F21").Value = 696
$ws4.Range("F22").Value = 609
$ws4.Range("F23").Value = 350
$ws4.Range("F25").Value = 74
$ws4.Range("F29").Value = 1253
$ws4.Range("F30").Value = 23
$ws4.Range("F35").Value = 1421
$ws4.Range("F38").Value = 351
$ws4.Range("F39").Value = 531
$ws4.Range("F41").Value = 367
$ws4.Range("F42").Value = 110
$ws4.Range("F46").Value = 650
$ws4.Range("F48").Value = 142
$ws4.Range("F49").Value = 128
